# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.967.46"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.879.98"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.96"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4997"
$ws.Range("E7").Value = "  -4.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09174"
$ws.Range("E9").Value = "  -5.62%  "
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.63"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.323"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.877.29"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.261"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.40"
$ws.Range("E18").Value = "  -3.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06622"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.92"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.181"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "28.026.09"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "2.096.76"
$ws.Range("E26").Value = "  -2.00%  "
$ws.Range("E27").Value = "  -5.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.49"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.73"
$ws.Range("E29").Value = "  -2.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.40"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1054"
$ws.Range("E32").Value = "  -2.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.574"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.595"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.313"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06542"
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02400"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2176"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("E39").Value = "  +8.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.203"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6392"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.52"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.923"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.26"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.294"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.984"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.201"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.08"
$ws.Range("E51").Value = "  -3.37%  "
